$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '66.279.03'
$ws.Range('E2').Value = '  -1.17%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.550.48'
$ws.Range('E3').Value = '  +1.08%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '608.62'
$ws.Range('E5').Value = '  -0.19%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '144.36'
$ws.Range('E6').Value = '  -2.69%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.548.06'
$ws.Range('E7').Value = '  +1.02%  '

$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('E9').Value = '  +0.32%  '

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.137'
$ws.Range('E10').Value = '  -3.94%  '

$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '8.07'
$ws.Range('E11').Value = '  -0.02%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.411'
$ws.Range('E12').Value = '  -2.69%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.146.13'
$ws.Range('E13').Value = '  +0.97%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000208'
$ws.Range('E14').Value = '  -3.84%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '30.11'
$ws.Range('E15').Value = '  -4.66%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.540.75'
$ws.Range('E16').Value = '  +0.69%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '66.342.40'
$ws.Range('E17').Value = '  -1.18%  '

$ws.Range('E18').Value = '  -0.87%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.93'
$ws.Range('E19').Value = '  +1.08%  '

$ws.Range('E20').Value = '  -2.62%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.92'
$ws.Range('E21').Value = '  -2.96%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '425.93'
$ws.Range('E22').Value = '  -2.51%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.603'
$ws.Range('E23').Value = '  -0.92%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '79.00'
$ws.Range('E24').Value = '  -0.91%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.689.87'
$ws.Range('E25').Value = '  +1.09%  '

$ws.Range('E26').Value = '  -0.01%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000121'
$ws.Range('E27').Value = '  +1.13%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.10'
$ws.Range('E28').Value = '  -1.59%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.17'
$ws.Range('E29').Value = '  -6.18%  '

$ws.Range('E30').Value = '  -1.44%  '

$ws.Range('E31').Value = '  +0.16%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.49'
$ws.Range('E32').Value = '  -5.97%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.160'
$ws.Range('E33').Value = '  -3.96%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '25.35'
$ws.Range('E34').Value = '  -0.91%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.539.23'
$ws.Range('E35').Value = '  +0.99%  '

$ws.Range('E36').Value = '  -0.01%  '

$ws.Range('E37').Value = '  -2.85%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '7.82'
$ws.Range('E38').Value = '  -2.60%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.63'
$ws.Range('E39').Value = '  -5.64%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.02%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '173.71'
$ws.Range('E41').Value = '  -1.37%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0857'
$ws.Range('E42').Value = '  -4.31%  '

$ws.Range('E43').Value = '  -2.80%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.894'
$ws.Range('E44').Value = '  -0.22%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.91'
$ws.Range('E45').Value = '  -6.47%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '45.64'
$ws.Range('E46').Value = '  -1.25%  '

$ws.Range('E47').Value = '  -1.42%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '26.00'
$ws.Range('E48').Value = '  -7.14%  '

$ws.Range('E49').Value = '  -1.65%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.14'
$ws.Range('E50').Value = '  -4.49%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.943'
$ws.Range('E51').Value = '  -5.33%  '
